$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2183.8
$ws.Range("I38").Value = 229.75
$ws.Range("K38").Value = 689.25
$ws.Range("M38").Value = -317.25

$ws.Range("H40").Value = 3472
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325

$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H86").Value = 1750
$ws.Range("I86").Value = 1750
$ws.Range("K86").Value = 1750
$ws.Range("M86").Value = -627

$ws.Range("H89").Value = 1750
$ws.Range("I89").Value = 1750
$ws.Range("K89").Value = 8750
$ws.Range("M89").Value = -3134

$ws.Range("H96").Value = 545.8333
$ws.Range("I96").Value = 447.16666
$ws.Range("J96").Value = 644.5
$ws.Range("K96").Value = 1341.49998
$ws.Range("L96").Value = 1933.5
$ws.Range("M96").Value = 31.50001999999995
$ws.Range("N96").Value = -4679.5

$ws.Range("H99").Value = 959.8
$ws.Range("I99").Value = 199.5
$ws.Range("J99").Value = 1466.6666
$ws.Range("K99").Value = 598.5
$ws.Range("L99").Value = 4399.9998
$ws.Range("M99").Value = 899.5
$ws.Range("N99").Value = -7395.9998

$ws.Range("H131").Value = 6779.6
$ws.Range("I131").Value = 3474.5
$ws.Range("K131").Value = 10423.5
$ws.Range("M131").Value = -5383.5

$ws.Range("H138").Value = 2360.889
$ws.Range("I138").Value = 1382.6666
$ws.Range("J138").Value = 2850
$ws.Range("K138").Value = 4147.9998
$ws.Range("L138").Value = 8550
$ws.Range("M138").Value = 992.0002000000004
$ws.Range("N138").Value = -18830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 199.57143
$ws.Range("I5").Value = 217.5
$ws.Range("J5").Value = 175.66667
$ws.Range("K5").Value = 217.5
$ws.Range("L5").Value = 175.66667
$ws.Range("M5").Value = -105.5
$ws.Range("N5").Value = -399.66667

$ws.Range("H45").Value = 4050
$ws.Range("I45").Value = 3900
$ws.Range("K45").Value = 3900
$ws.Range("M45").Value = -3523

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H132").Value = 4671.125
$ws.Range("I132").Value = 2661.5
$ws.Range("K132").Value = 7984.5
$ws.Range("M132").Value = -5454.5

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 199.57143
$ws.Range("I4").Value = 217.5
$ws.Range("J4").Value = 175.66667
$ws.Range("K4").Value = 217.5
$ws.Range("L4").Value = 175.66667
$ws.Range("M4").Value = -102.5
$ws.Range("N4").Value = -405.66667

$ws.Range("H105").Value = 34398.832
$ws.Range("I105").Value = 1098.75
$ws.Range("J105").Value = 100999
$ws.Range("K105").Value = 1098.75
$ws.Range("L105").Value = 100999
$ws.Range("M105").Value = 648.25
$ws.Range("N105").Value = -104493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 642.4
$ws.Range("I22").Value = 443.125
$ws.Range("K22").Value = 443.125
$ws.Range("M22").Value = -93.125

$ws.Range("H36").Value = 6500
$ws.Range("J36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("N36").Value = -10776

$ws.Range("H40").Value = 6500
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10320

$ws.Range("H99").Value = 5605.8667
$ws.Range("I99").Value = 5010.5713
$ws.Range("K99").Value = 5010.5713
$ws.Range("M99").Value = -3512.5713

$ws.Range("H126").Value = 5605.8667
$ws.Range("I126").Value = 5010.5713
$ws.Range("K126").Value = 15031.7139
$ws.Range("M126").Value = -12561.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1081.8
$ws.Range("I113").Value = 636.3333
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1908.9999
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = 261.0001
$ws.Range("N113").Value = -9590

$ws.Range("H117").Value = 6456.4
$ws.Range("I117").Value = 543.2
$ws.Range("J117").Value = 12369.6
$ws.Range("K117").Value = 1629.6
$ws.Range("L117").Value = 37108.8
$ws.Range("M117").Value = 1812.4
$ws.Range("N117").Value = -43992.8

$ws.Range("H134").Value = 2866.5
$ws.Range("I134").Value = 1599.75
$ws.Range("K134").Value = 4799.25
$ws.Range("M134").Value = 270.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 198
$ws.Range("I22").Value = 199
$ws.Range("K22").Value = 199
$ws.Range("M22").Value = 330

$ws.Range("H102").Value = 1137.3334
$ws.Range("I102").Value = 1137.3334
$ws.Range("K102").Value = 1137.3334
$ws.Range("M102").Value = 484.6666

$ws.Range("H107").Value = 437
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 504.4
$ws.Range("K107").Value = 100
$ws.Range("L107").Value = 504.4
$ws.Range("M107").Value = 1820
$ws.Range("N107").Value = -4344.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 771.9091
$ws.Range("I22").Value = 578.4
$ws.Range("K22").Value = 578.4
$ws.Range("M22").Value = -283.4

$ws.Range("H27").Value = 771.9091
$ws.Range("I27").Value = 578.4
$ws.Range("K27").Value = 578.4
$ws.Range("M27").Value = -471.4

$ws.Range("H46").Value = 345
$ws.Range("I46").Value = 345
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 345
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -157
$ws.Range("N46").ClearContents()

$ws.Range("H103").Value = 30000.666
$ws.Range("J103").Value = 30000.666
$ws.Range("L103").Value = 30000.666
$ws.Range("N103").Value = -32344.666

$ws.Range("H123").Value = 15000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 504
$ws.Range("I14").Value = 504
$ws.Range("K14").Value = 504
$ws.Range("M14").Value = -336

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H69").Value = 25987.25
$ws.Range("J69").Value = 25987.25
$ws.Range("L69").Value = 25987.25
$ws.Range("N69").Value = -27485.25

$ws.Range("H72").Value = 25987.25
$ws.Range("J72").Value = 25987.25
$ws.Range("L72").Value = 77961.75
$ws.Range("N72").Value = -85449.75

$ws.Range("H113").Value = 2994.4
$ws.Range("I113").Value = 1836
$ws.Range("J113").Value = 3766.6667
$ws.Range("K113").Value = 5508
$ws.Range("L113").Value = 11300.0001
$ws.Range("M113").Value = -3338
$ws.Range("N113").Value = -15640.0001

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 1298.8182
$ws.Range("J132").Value = 900
$ws.Range("L132").Value = 2700
$ws.Range("N132").Value = -7760
